$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B8").NumberFormat = "@"

$ws.Range("B2").Value = "0.007148449489526548"
$ws.Range("B3").Value = "19.0"
$ws.Range("B4").Value = "0.0006320705701557373"
$ws.Range("B5").Value = "5.179405571978401"
$ws.Range("B6").Value = "0.00016045503652754083"
$ws.Range("B7").Value = "59.0"
$ws.Range("B8").Value = "9.5"

$ws.Range("B2:B8").ClearFormats()
